$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 5.850999999999996
$ws.Range("A3").Value = -21.36680000000002
$ws.Range("B5").Value = 5.154600000000004
$ws.Range("D5").Value = -8.705099999999998
$ws.Range("D9").Value = -8.5336
$ws.Range("D11").Value = -8.293800000000003
$ws.Range("A14").Value = -20.51159999999998
$ws.Range("A21").Value = -21.37400000000001
$ws.Range("D21").Value = -8.150500000000003
$ws.Range("A23").Value = -21.33140000000002
$ws.Range("A25").Value = -22.42640000000003
